# Applies the resume content edit described by the commit diff:
#  - Rewrites the Skills bullet list with new soft-skill bullets
#  - Reworks the Experience section (Undergraduate Research Assistant bullets,
#    and renames "Executive Team Member" -> "Summer Intern - Sharpen Up
#    Internship Program (Rotational)" with new org/dates/bullets)
#  - Reworks the Projects section (renames the two existing projects with new
#    titles/dates/tech-stacks/bullets) and removes the trailing
#    "Full-stack Financial Assistant | Hackathon Project" entry entirely.
#
# Text substitutions are done via Find.Execute (search only, no built-in
# Replace) followed by a direct Range.Text assignment. Using Find's own
# Replace argument silently "smart-quotes" straight apostrophes in the
# replacement string, which would corrupt "OpenAI's" below, so we avoid it.

$d = $word.ActiveDocument

function Replace-Text($find, $replace, [bool]$wholeWord = $false) {
    $rng = $d.Content
    $ok = $rng.Find.Execute($find, $true, $wholeWord, $false, $false, $false, $true, 1, $false)
    if (-not $ok) {
        throw "Find failed for: $find"
    }
    $rng.Text = $replace
}

# --- Skills section -------------------------------------------------------
Replace-Text "• Programming languages: Java, Python, JavaScript, C, C++, HTML, CSS, React, Node.JS, Express.JS, SQL" "• Cloud computing and AI integration"
Replace-Text "• Software tools: VS Code, Git, Github, Gitlab, Unity, Unreal Engine" "• Data-driven decision making"
Replace-Text "• Cloud and Data Management: AWS, Azure, GCP basics" "• Cross-functional team collaboration"
Replace-Text "• AI/ML: Familiar with basic AI/ML concepts and integration" "• Adapting to new technologies and tools"
Replace-Text "• Algorithm and Data Structures: Experienced with university-level study" "• Detail-oriented analysis and problem solving"
Replace-Text "• Professional Skills: Adaptable, Excellent communication, Detail-oriented, Leadership, Time Management" "• Strong communication and interpersonal skills"

# --- Experience: Undergraduate Research Assistant -------------------------
Replace-Text "❖ Undergraduate Research Assistant (Node, React, JS)" "❖ Undergraduate Research Assistant"
Replace-Text "University of Calgary, Calgary, AB" "University of Calgary – Calgary, AB"
Replace-Text "• Developed an automated workflow using Node and React for extracting detailed data insights in a timely manner." "• Developed automated workflows combining cutting-edge technologies like Node, React, and OpenAI's Whisper to enhance data processing efficiency."
Replace-Text "• Collected and processed multi-modal data (videos, spoken recordings, biometric data) for research in information needs." "• Collaborated cross-functionally to create data visualizations for large datasets, honing skills in data-driven insights."
Replace-Text "• Adapted quickly to new tools and technologies to enhance research data analysis processes, showcasing adaptability." "• Demonstrated adaptability by mastering new technology stacks and applied problem-solving to automate data analysis workflows."

# --- Experience: Executive Team Member -> Summer Intern --------------------
Replace-Text "❖ Executive Team Member" "❖ Summer Intern – Sharpen Up Internship Program (Rotational)"
Replace-Text "Dec 2021 - Apr 2022" "Jun 2025 – Aug 2025"
Replace-Text "Model United Nations (MUN) at Dar Jana International School" "Viatris Egypt – Cairo, Egypt"
Replace-Text "• Organized and prepared event documents, ensuring seamless execution of MUN conferences." "• Gained a broad understanding of pharmaceutical operations by rotating across various departments, contributing to cross-functional initiatives."
Replace-Text "• Managed participant engagement and facilitated communication, enhancing collaborative problem-solving skills." "• Assisted in developing dashboards using Excel and SQL for real-time supply chain performance monitoring."
Replace-Text "• Demonstrated leadership by acting as a spokesman, guiding event procedures effectively." "• Learned and applied pharmaceutical compliance and process improvement principles, enhancing analytical and organizational skills."

# --- Projects: Self-Checkout Machine Software -> SceneBook -----------------
# (Must run before the standalone "Java" subtitle replacement below, since
# this heading text also contains the substring "Java".)
Replace-Text "❖ Self-Checkout Machine Software (Java)" "❖ SceneBook: Multi-Theatre Aggregation & Ticketing Platform"
Replace-Text "Sep 2023 – Dec 2023" "Jan 2025 – Apr 2025"

# Standalone "Java" tech-stack line: match whole word only, and only after
# the two other "Java"-containing strings above have already been replaced,
# so this now uniquely targets the subtitle paragraph.
Replace-Text "Java" "PostgreSQL, Node.js, React" $true

Replace-Text "• Collaborated with a team of 20 to design and develop user-friendly software for self-checkout systems." "• Built a full-stack platform enhancing user experience via unified movie listings, showtimes, and booking flows."
Replace-Text "• Focused on efficient transaction handling and integrated real-world use case functionalities." "• Designed a scalable PostgreSQL database schema, facilitating multi-theatre management and seamless admin functionality."
Replace-Text "• Contributed to interface design improvements that enhanced customer experience." "• Engineered intuitive frontend features with React, achieving robust connectivity and user engagement."

# --- Projects: Educational Assessment Web App -> EventEcho -----------------
Replace-Text "❖ Educational Assessment Web App (JS, CSS, HTML)" "❖ EventEcho – Full-Stack Event Management Web Application"
Replace-Text "Jan 2024 – Apr 2024" "Sep 2024 – Dec 2024"
Replace-Text "JavaScript, CSS, HTML" "Node.js, React, JWT, PostgreSQL, Docker"
Replace-Text "• Developed a dynamic web application with a team of 5 for creating educational assessments with real-time feedback." "• Implemented token-based authentication securing user login and admin access with least-privilege principles."
Replace-Text "• Implemented features to randomly generate questions and provide immediate grading." "• Developed RESTful APIs supporting event registration, user management, and real-time data transactions."
Replace-Text "• Prioritized user-friendly navigation and engaging test experiences." "• Integrated cloud services for scalable event hosting and user data management."

# --- Projects: remove the "Full-stack Financial Assistant | Hackathon
#     Project" entry entirely (heading + subtitle + 3 bullet paragraphs).
$startPara = $null
$endPara = $null
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs($i)
    $text = $p.Range.Text
    if ($text -like "*Full-stack Financial Assistant | Hackathon Project*") {
        $startPara = $i
    }
    if ($startPara -ne $null -and $text -like "*Demonstrated strong team collaboration and project management under time constraints.*") {
        $endPara = $i
        break
    }
}

if ($startPara -eq $null -or $endPara -eq $null) {
    throw "Could not locate the Hackathon Project paragraph block to delete."
}

$rangeStart = $d.Paragraphs($startPara).Range.Start
$rangeEnd = $d.Paragraphs($endPara).Range.End
$d.Range($rangeStart, $rangeEnd).Delete()

Write-Output "Done."
